$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Fill in all the numeric/date/time values first (these don't touch the
# shared-strings table, so order among themselves does not matter).
# ---------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = 44117
$ws.Range("B2").Value = 0.44791666666666669
$ws.Range("C2").Value = 0.4513888888888889
$ws.Range("E2").Value = 0.05

# Row 3
$ws.Range("A3").Value = 44117
$ws.Range("B3").Value = 0.4513888888888889
$ws.Range("C3").Value = 0.4548611111111111
$ws.Range("E3").Value = 0.05

# Row 4
$ws.Range("D4").Value = 0.3

# Row 5
$ws.Range("A5").Value = 44117
$ws.Range("B5").Value = 0.47569444444444442
$ws.Range("D5").Value = 0.50347222222222221
$ws.Range("E5").Value = 0.4

# Row 7
$ws.Range("A7").Value = 44124
$ws.Range("B7").Value = 0.4284722222222222
$ws.Range("C7").Value = 0.43055555555555558
$ws.Range("E7").Value = 0.03

# Row 8
$ws.Range("A8").Value = 44124
$ws.Range("B8").Value = 0.43124999999999997
$ws.Range("C8").Value = 0.46180555555555558
$ws.Range("E8").Value = 0.44

# Row 9
$ws.Range("D9").Value = 0.2

# Row 10
$ws.Range("B10").Value = 0.47569444444444442
$ws.Range("C10").Value = 0.4826388888888889
$ws.Range("E10").Value = 0.1

# Row 12
$ws.Range("A12").Value = 44131
$ws.Range("B12").Value = 0.52083333333333337
$ws.Range("C12").Value = 0.52430555555555558
$ws.Range("E12").Value = 0.05

# Row 13
$ws.Range("B13").Value = 0.52430555555555558
$ws.Range("C13").Value = 0.54861111111111105
$ws.Range("E13").Value = 0.35

# Row 15
$ws.Range("A15").Value = 44138
$ws.Range("B15").Value = 0.42708333333333331
$ws.Range("C15").Value = 0.43402777777777773
$ws.Range("E15").Value = 0.1

# Row 16
$ws.Range("B16").Value = 0.43402777777777773
$ws.Range("C16").Value = 0.46527777777777773
$ws.Range("E16").Value = 0.45

# ---------------------------------------------------------------------
# Now fill in the "description" (column F) text values, in the same
# order the original author typed them (session rows first, then the
# "Break" label backfilled afterwards) so the shared-strings table ends
# up with entries registered in the right order.
# ---------------------------------------------------------------------

$ws.Range("F2").Value = "Making scene"
$ws.Range("F3").Value = "Adding player and ridgidbody and starting script "
$ws.Range("F5").Value = "making player script"
$ws.Range("F4").Value = "Break"
$ws.Range("F7").Value = "Making scene"
$ws.Range("F8").Value = "Making movement code and camera follow"
$ws.Range("F10").Value = "Adding scripts "
$ws.Range("F12").Value = "Making scene"
$ws.Range("F13").Value = "Coding Jump"
$ws.Range("F15").Value = "Setting scene up "
$ws.Range("F16").Value = "Adding endless movement and respawn, adding to player and setting movement speed. "

# ---------------------------------------------------------------------
# Apply number formats. The first cell of each kind defines the style;
# the rest are formatted by copying that cell's format only, so every
# cell of the same kind shares a single reused cell style (just like
# Excel's own format-painter / fill-down behaviour).
# ---------------------------------------------------------------------

# Date style (numFmtId 14, "mm-dd-yy") applied to column A entries
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$dateCells = @("A3", "A5", "A7", "A8", "A12", "A15")
foreach ($cell in $dateCells) {
    $ws.Range($cell).PasteSpecial(-4122)
}

# Time style (numFmtId 20, "h:mm") applied to the Start/End/break columns
$ws.Range("B2").NumberFormat = "h:mm"
$ws.Range("B2").Copy()
$timeCells = @("C2", "B3", "C3", "B5", "D5", "B7", "C7", "B8", "C8", "B10", "C10", "B12", "C12", "B13", "C13", "B15", "C15", "B16", "C16")
foreach ($cell in $timeCells) {
    $ws.Range($cell).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Move the active selection to F17, matching the final cursor position after data entry
$ws.Range("F17").Select()
